$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 74 (the "私の未来" post), shifting all subsequent
# rows up by one. This matches the target edit where that post was removed
# from the sheet and everything below it renumbered accordingly.
$ws.Rows.Item(74).Delete()
